# Daily attendance processing - move "System" to the front of the
# "Recorded By" (column G) comma-separated list of recorders, on every
# row of the active sheet where "System" (exact case) is present but
# not already the first entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $text = $cell.Value2

    if ($text -ne $null -and $text -ne "") {
        $parts = $text -split ",\s*"

        if ($parts.Count -gt 1) {
            $hasSystem = $false
            foreach ($p in $parts) {
                if ($p.Equals("System")) { $hasSystem = $true }
            }
            $firstIsSystem = $parts[0].Equals("System")

            if ($hasSystem -and (-not $firstIsSystem)) {
                $newParts = @("System")
                foreach ($p in $parts) {
                    if (-not $p.Equals("System")) {
                        $newParts += $p
                    }
                }
                $cell.Value2 = [string]::Join(", ", $newParts)
            }
        }
    }
}
